{"js": "// Remove the trailing \"Ver no Jupiter ...\" / copyright paragraphs (plus the\n// blank paragraph separating them from the requirements line) that used to\n// follow \"LOQ4233: Gest\u00e3o de Neg\u00f3cios (Requisito fraco)\".\n\nconst body = context.document.body;\n\n// Anchor on the paragraph that must remain (end of the Requisitos section).\nconst results = body.search(\"LOQ4233: Gest\u00e3o de Neg\u00f3cios (Requisito fraco)\", {\n  matchCase: true,\n});\nresults.load(\"items\");\nawait context.sync();\n\nconst anchorRange = results.items[0];\nconst anchorParas = anchorRange.paragraphs;\nanchorParas.load(\"items\");\nawait context.sync();\n\nconst anchorPara = anchorParas.items[0];\n\n// The three paragraphs to delete: the blank line, the \"Ver no Jupiter...\"\n// line, and the \"\u00a9 2020 ...\" copyright line.\nconst firstToDelete = anchorPara.getNext();\nconst secondToDelete = firstToDelete.getNext();\nconst thirdToDelete = secondToDelete.getNext();\n\nfirstToDelete.delete();\nsecondToDelete.delete();\nthirdToDelete.delete();\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" / copyright paragraphs (plus the\n# blank paragraph separating them from the requirements line) that used to\n# follow \"LOQ4233: Gest\u00e3o de Neg\u00f3cios (Requisito fraco)\".\n\n$d = $word.ActiveDocument\n\n# Anchor on the paragraph that must remain (end of the Requisitos section).\n$anchor = $d.Content\n$null = $anchor.Find.Execute(\"LOQ4233: Gest\u00e3o de Neg\u00f3cios (Requisito fraco)\")\n\n$anchorPara = $anchor.Paragraphs(1)\n\n# The three paragraphs to delete: the blank line, the \"Ver no Jupiter...\"\n# line, and the \"\u00a9 2020 ...\" copyright line.\n$firstToDelete = $anchorPara.Next()\n$secondToDelete = $firstToDelete.Next()\n$thirdToDelete = $secondToDelete.Next()\n\n$deleteRange = $d.Range($firstToDelete.Range.Start, $thirdToDelete.Range.End)\n$deleteRange.Delete()\n"}
